$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking stat values stay stored as text (matching original t="str" cells)
$ws.Range("C2:F9").NumberFormat = "@"

# New row 9 needs the player/team labels copied down as well (reuse exact text from row 8)
$ws.Range("A9").Value = $ws.Range("A8").Text
$ws.Range("B9").Value = $ws.Range("B8").Text

# Row 2: runs, balls, fours, sixes
$ws.Range("C2").Value = "50"
$ws.Range("D2").Value = "26"
$ws.Range("E2").Value = "6"
$ws.Range("F2").Value = "3"

# Row 3
$ws.Range("C3").Value = "30"
$ws.Range("D3").Value = "32"
$ws.Range("E3").Value = "2"
$ws.Range("F3").Value = "0"

# Row 4
$ws.Range("C4").Value = "107"
$ws.Range("D4").Value = "60"
$ws.Range("E4").Value = "14"
$ws.Range("F4").Value = "3"

# Row 5
$ws.Range("C5").Value = "18"
$ws.Range("D5").Value = "11"
$ws.Range("E5").Value = "2"
$ws.Range("F5").Value = "1"

# Row 6
$ws.Range("C6").Value = "19"
$ws.Range("D6").Value = "11"
$ws.Range("E6").Value = "3"
$ws.Range("F6").Value = "0"

# Row 7
$ws.Range("C7").Value = "41"
$ws.Range("D7").Value = "35"
$ws.Range("E7").Value = "6"
$ws.Range("F7").Value = "0"

# Row 8
$ws.Range("C8").Value = "5"
$ws.Range("D8").Value = "6"
$ws.Range("E8").Value = "1"
$ws.Range("F8").Value = "0"

# Row 9 (newly added record)
$ws.Range("C9").Value = "15"
$ws.Range("D9").Value = "19"
$ws.Range("E9").Value = "2"
$ws.Range("F9").Value = "0"
